$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new Wins/Losses/Ties columns (AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, border, centered) from an existing header cell (AC1) to the new headers
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill in team record values for data rows 2-59
for ($r = 2; $r -le 59; $r++) {
    $ws.Cells.Item($r, 30).Value = 73   # AD = column 30
    $ws.Cells.Item($r, 31).Value = 89   # AE = column 31
    $ws.Cells.Item($r, 32).Value = 0    # AF = column 32
}
